$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = '이닉스(구.이닉스정호)'
$ws.Cells.Item(2, 2).Value = '2024.01.11~01.17'
$ws.Cells.Item(2, 3).Value = '9,200~11,000'
$ws.Cells.Item(2, 4).Value = '-'
$ws.Cells.Item(2, 5).Value = 27600
$ws.Cells.Item(2, 6).Value = '삼성증권'

$ws.Cells.Item(3, 1).Value = '에이치비인베스트먼트'
$ws.Cells.Item(3, 2).Value = '2024.01.08~01.12'
$ws.Cells.Item(3, 3).Value = '2,400~2,800'
$ws.Cells.Item(3, 4).Value = '-'
$ws.Cells.Item(3, 5).Value = 16000
$ws.Cells.Item(3, 6).Value = 'NH투자증권'

$ws.Cells.Item(4, 1).Value = '현대힘스'
$ws.Cells.Item(4, 2).Value = '2024.01.08~01.12'
$ws.Cells.Item(4, 3).Value = '5,000~6,300'
$ws.Cells.Item(4, 4).Value = '-'
$ws.Cells.Item(4, 5).Value = 43535
$ws.Cells.Item(4, 6).Value = '미래에셋증권'

$ws.Cells.Item(5, 1).Value = '포스뱅크'
$ws.Cells.Item(5, 2).Value = '2024.01.05~01.11'
$ws.Cells.Item(5, 3).Value = '13,000~15,000'
$ws.Cells.Item(5, 4).Value = '-'
$ws.Cells.Item(5, 5).Value = 19500
$ws.Cells.Item(5, 6).Value = '하나증권'

$ws.Cells.Item(6, 1).Value = '하나스팩30호'
$ws.Cells.Item(6, 2).Value = '2023.12.08~12.11'
$ws.Cells.Item(6, 3).Value = '2,000~2,000'
$ws.Cells.Item(6, 4).Value = '-'
$ws.Cells.Item(6, 5).Value = 14000
$ws.Cells.Item(6, 6).Value = '하나증권'

$ws.Cells.Item(7, 1).Value = '디에스단석(구,단석산업)'
$ws.Cells.Item(7, 2).Value = '2023.12.05~12.11'
$ws.Cells.Item(7, 3).Value = '79,000~89,000'
$ws.Cells.Item(7, 4).Value = '-'
$ws.Cells.Item(7, 5).Value = 96380
$ws.Cells.Item(7, 6).Value = 'KB증권,NH투자증권'

$ws.Cells.Item(8, 1).Value = 'IBKS스팩23호'
$ws.Cells.Item(8, 2).Value = '2023.12.04~12.08'
$ws.Cells.Item(8, 3).Value = '2,000~2,000'
$ws.Cells.Item(8, 4).Value = '2000'
$ws.Cells.Item(8, 5).Value = 8000
$ws.Cells.Item(8, 6).Value = '아이비케이투자증권'

$ws.Cells.Item(9, 1).Value = '블루엠텍'
$ws.Cells.Item(9, 2).Value = '2023.11.22~11.28'
$ws.Cells.Item(9, 3).Value = '15,000~19,000'
$ws.Cells.Item(9, 4).Value = '19000'
$ws.Cells.Item(9, 5).Value = 21000
$ws.Cells.Item(9, 6).Value = '하나증권,키움증권'

$ws.Cells.Item(10, 1).Value = 'LS머트리얼즈'
$ws.Cells.Item(10, 2).Value = '2023.11.22~11.28'
$ws.Cells.Item(10, 3).Value = '4,400~5,500'
$ws.Cells.Item(10, 4).Value = '6000'
$ws.Cells.Item(10, 5).Value = 64350
$ws.Cells.Item(10, 6).Value = '키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권'

$ws.Cells.Item(11, 1).Value = '삼성스팩9호'
$ws.Cells.Item(11, 2).Value = '2023.11.20~11.21'
$ws.Cells.Item(11, 3).Value = '2,000~2,000'
$ws.Cells.Item(11, 4).Value = '2000'
$ws.Cells.Item(11, 5).Value = 20000
$ws.Cells.Item(11, 6).Value = '삼성증권'

$ws.Cells.Item(12, 1).Value = '교보스팩15호'
$ws.Cells.Item(12, 2).Value = '2023.11.20~11.21'
$ws.Cells.Item(12, 3).Value = '2,000~2,000'
$ws.Cells.Item(12, 4).Value = '2000'
$ws.Cells.Item(12, 5).Value = 7000
$ws.Cells.Item(12, 6).Value = '교보증권'

$ws.Cells.Item(13, 1).Value = '케이엔에스'
$ws.Cells.Item(13, 2).Value = '2023.11.16~11.22'
$ws.Cells.Item(13, 3).Value = '19,000~22,000'
$ws.Cells.Item(13, 4).Value = '23000'
$ws.Cells.Item(13, 5).Value = 14250
$ws.Cells.Item(13, 6).Value = '신영증권'

$ws.Cells.Item(14, 1).Value = 'NH스팩30호'
$ws.Cells.Item(14, 2).Value = '2023.11.15~11.16'
$ws.Cells.Item(14, 3).Value = '2,000~2,000'
$ws.Cells.Item(14, 4).Value = '2000'
$ws.Cells.Item(14, 5).Value = 16000
$ws.Cells.Item(14, 6).Value = 'NH투자증권'

$ws.Cells.Item(15, 1).Value = '와이바이오로직스'
$ws.Cells.Item(15, 2).Value = '2023.11.10~11.16'
$ws.Cells.Item(15, 3).Value = '9,000~11,000'
$ws.Cells.Item(15, 4).Value = '9000'
$ws.Cells.Item(15, 5).Value = 13500
$ws.Cells.Item(15, 6).Value = '유안타증권'

$ws.Cells.Item(16, 1).Value = '에이텀'
$ws.Cells.Item(16, 2).Value = '2023.11.09~11.15'
$ws.Cells.Item(16, 3).Value = '23,000~30,000'
$ws.Cells.Item(16, 4).Value = '18000'
$ws.Cells.Item(16, 5).Value = 14950
$ws.Cells.Item(16, 6).Value = '하나증권'

$ws.Cells.Item(17, 1).Value = '에이에스텍'
$ws.Cells.Item(17, 2).Value = '2023.11.07~11.13'
$ws.Cells.Item(17, 3).Value = '21,000~25,000'
$ws.Cells.Item(17, 4).Value = '28000'
$ws.Cells.Item(17, 5).Value = 29547
$ws.Cells.Item(17, 6).Value = '미래에셋증권'

$ws.Cells.Item(18, 1).Value = '그린리소스'
$ws.Cells.Item(18, 2).Value = '2023.11.03~11.09'
$ws.Cells.Item(18, 3).Value = '11,000~14,000'
$ws.Cells.Item(18, 4).Value = '17000'
$ws.Cells.Item(18, 5).Value = 18040
$ws.Cells.Item(18, 6).Value = 'NH투자증권'

$ws.Cells.Item(19, 1).Value = '한선엔지니어링'
$ws.Cells.Item(19, 2).Value = '2023.11.02~11.08'
$ws.Cells.Item(19, 3).Value = '5,200~6,000'
$ws.Cells.Item(19, 4).Value = '7000'
$ws.Cells.Item(19, 5).Value = 22100
$ws.Cells.Item(19, 6).Value = '대신증권'

$ws.Cells.Item(20, 1).Value = '에코아이'
$ws.Cells.Item(20, 2).Value = '2023.11.01~11.07'
$ws.Cells.Item(20, 3).Value = '28,500~34,700'
$ws.Cells.Item(20, 4).Value = '34700'
$ws.Cells.Item(20, 5).Value = 59251
$ws.Cells.Item(20, 6).Value = 'KB증권'

$ws.Cells.Item(21, 1).Value = '동인기연(유가)'
$ws.Cells.Item(21, 2).Value = '2023.11.01~11.07'
$ws.Cells.Item(21, 3).Value = '33,000~37,000'
$ws.Cells.Item(21, 4).Value = '30000'
$ws.Cells.Item(21, 5).Value = 60654
$ws.Cells.Item(21, 6).Value = 'NH투자증권'
